$wb = $excel.ActiveWorkbook

# Sheet handles (by position, 1-based):
# 1 = c CO2 Capture, 2 = v CO2 Capture, 3 = c CO2 Compression,
# 4 = v CO2 Compression, 5 = c Power Station, 6 = v Power Station
$wsVCapture     = $wb.Worksheets.Item(2)
$wsVCompression = $wb.Worksheets.Item(4)
$wsVPower       = $wb.Worksheets.Item(6)

# --- v Power Station ("sheet6"): rename scenarios and add the new EU-BECCS / charcoal row ---
$wsVPower.Range("A4").Value = "EU-2010"
$wsVPower.Range("A5").Value = "EU-1990"
$wsVPower.Range("A6").Value = "EU-2000"

$wsVPower.Range("A7").Value = "EU-BECCS"
$wsVPower.Range("B7").Value = 0.85
$wsVPower.Range("C7").Value = 0.85
$wsVPower.Range("D7").Value = "charcoal"

# --- v CO2 Capture ("sheet2"): rename scenarios, update EU-1990 row, add EU-BECCS row ---
$wsVCapture.Range("A4").Value = "EU-2010"
$wsVCapture.Range("A5").Value = "EU-2000"
$wsVCapture.Range("A6").Value = "EU-1990"
$wsVCapture.Range("B6").Value = 0.4
$wsVCapture.Range("C6").Value = 2

$wsVCapture.Range("A7").Value = "EU-BECCS"
$wsVCapture.Range("B7").Value = 0.8
$wsVCapture.Range("C7").Value = 1.2

# --- v CO2 Compression ("sheet4"): rename scenarios, add EU-BECCS row ---
$wsVCompression.Range("A4").Value = "EU-2010"
$wsVCompression.Range("A5").Value = "EU-2000"
$wsVCompression.Range("A6").Value = "EU-1990"

$wsVCompression.Range("A7").Value = "EU-BECCS"
$wsVCompression.Range("B7").Value = 0.25
$wsVCompression.Range("C7").Value = 0.02

# --- View / selection state ---
[void]$wsVCompression.Range("A4:A7").Select()
[void]$wsVPower.Range("L20").Select()

# Make "v CO2 Capture" the active sheet/tab, with E16 selected
[void]$wsVCapture.Activate()
[void]$wsVCapture.Range("E16").Select()
